# Add a new "LoginCenter" scene row (row 13) to the StartSceneConfig sheet,
# mirroring the existing "Account" row (row 12) for Id=8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 12's formatting onto row 13 first (Id/Zone/SceneType/OuterPort
# columns all reuse the look of the row above it).
$ws.Range("C12:H12").Copy()
$ws.Range("C13:H13").PasteSpecial(-4122)

# New scene config values.
$ws.Range("C13").Value = 8
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = "LoginCenter"
$ws.Range("G13").Value = "LoginCenter"
$ws.Range("H13").ClearContents()

# Match font/alignment used by the "Account"/"LoginCenter" style column.
$ws.Range("F13:G13").Font.Name = "微软雅黑"
$ws.Range("F13:G13").Font.Size = 9
$ws.Range("F13:G13").HorizontalAlignment = -4131

$ws.Range("H13").Select()
